$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'65.812.91"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.38%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.166.96"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -5.04%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.01%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'571.47"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.81%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'171.42"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -4.17%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.01%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.595"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -3.19%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'3.162.27"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -5.02%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.124"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -3.51%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'6.59"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -4.05%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.392"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -4.09%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'3.716.16"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -5.01%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'0.135"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.50%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'27.36"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -4.72%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'65.753.86"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.29%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  -2.73%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'3.172.35"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -4.73%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'5.71"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.20%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'12.89"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -4.15%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'360.48"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -0.71%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'7.27"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -2.21%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.07%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'69.13"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -3.27%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').Value = "'WrappedeETH"
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = "'3.309.12"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -4.92%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value = "'Polygon"
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = "'0.494"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -4.98%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'0.0000114"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -6.94%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'9.85"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +3.01%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.57%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.08%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.09%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  -1.91%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = "'EthereumClassic"
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = "'22.08"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -3.57%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = "'NEARProtocol"
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = "'5.36"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -4.49%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.47%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'6.61"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -3.09%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'159.20"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -1.04%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'1.45"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -2.64%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.835"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -1.00%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +3.49%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'26.42"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -3.48%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'2.48"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -3.59%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'2.640.74"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -2.64%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'6.17"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -1.14%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'4.18"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -2.15%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'39.70"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.17%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.0658"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -1.47%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'327.96"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.58%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'23.94"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.33%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.0274"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -1.87%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.101"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -1.85%  "
$ws.Range('E51').Style = 'Normal'
